# Update odds values in "Jogos_da_Semana_FlashScore_2025-05-20" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("N2").Value = 2.01
$ws.Range("O2").Value = 1.89

# Row 3
$ws.Range("L3").Value = 1.14
$ws.Range("M3").Value = 5.5
$ws.Range("AB3").Value = 15

# Row 25
$ws.Range("G25").Value = 2.1
$ws.Range("I25").Value = 3.1
$ws.Range("J25").Value = 1.02
$ws.Range("K25").Value = 12
$ws.Range("T25").Value = 8.5
$ws.Range("W25").Value = 19

# Row 28
$ws.Range("G28").Value = 1.7
$ws.Range("H28").Value = 3.85
$ws.Range("I28").Value = 4.6
$ws.Range("K28").Value = 9
$ws.Range("L28").Value = 1.23
$ws.Range("M28").Value = 3.9
$ws.Range("N28").Value = 1.7
$ws.Range("O28").Value = 2.07
$ws.Range("P28").Value = 1.36
$ws.Range("Q28").Value = 3
$ws.Range("R28").Value = 1.7
$ws.Range("S28").Value = 2.05
$ws.Range("T28").Value = 7.7
$ws.Range("U28").Value = 9.25
$ws.Range("X28").Value = 13.5
$ws.Range("Y28").Value = 24
$ws.Range("Z28").Value = 9
$ws.Range("AA28").Value = 7.9
$ws.Range("AB28").Value = 15.5
$ws.Range("AC28").Value = 65
$ws.Range("AD28").Value = 450
$ws.Range("AI28").Value = 40
